{"js": "// Update the date heading and the multiplication answers in the practice table.\nconst body = context.document.body;\n\nconst pairs = [\n  [\"2023-08-06 Sunday\", \"2023-08-07 Monday\"],\n  [\"95\u00d714=1330\", \"47\u00d763=2961\"],\n  [\"63\u00d796=6048\", \"22\u00d783=1826\"],\n  [\"63\u00d778=4914\", \"91\u00d746=4186\"],\n  [\"49\u00d758=2842\", \"72\u00d792=6624\"],\n  [\"90\u00d778=7020\", \"32\u00d732=1024\"],\n  [\"43\u00d728=1204\", \"38\u00d711=418\"],\n  [\"47\u00d786=4042\", \"96\u00d729=2784\"],\n  [\"80\u00d774=5920\", \"16\u00d733=528\"],\n  [\"79\u00d742=3318\", \"22\u00d766=1452\"],\n  [\"33\u00d731=1023\", \"76\u00d761=4636\"],\n  [\"30\u00d777=2310\", \"92\u00d718=1656\"],\n  [\"65\u00d731=2015\", \"75\u00d764=4800\"],\n  [\"15\u00d716=240\", \"14\u00d783=1162\"],\n  [\"13\u00d739=507\", \"85\u00d728=2380\"],\n  [\"12\u00d733=396\", \"41\u00d755=2255\"],\n  [\"61\u00d753=3233\", \"58\u00d717=986\"],\n  [\"88\u00d720=1760\", \"72\u00d765=4680\"],\n  [\"65\u00d736=2340\", \"94\u00d734=3196\"],\n  [\"75\u00d771=5325\", \"28\u00d783=2324\"],\n  [\"92\u00d725=2300\", \"34\u00d795=3230\"],\n  [\"24\u00d732=768\", \"94\u00d754=5076\"],\n  [\"37\u00d796=3552\", \"31\u00d754=1674\"],\n  [\"74\u00d750=3700\", \"94\u00d793=8742\"],\n  [\"89\u00d779=7031\", \"76\u00d729=2204\"],\n  [\"44\u00d717=748\", \"80\u00d722=1760\"]\n];\n\n// Kick off a search for every old value, then resolve them all in one sync.\nconst searchResults = pairs.map(([oldText]) => body.search(oldText, { matchCase: true, matchWholeWord: false }));\nawait context.sync();\n\n// Replace every match found for each search with its corresponding new value.\nfor (let i = 0; i < pairs.length; i++) {\n  const [, newText] = pairs[i];\n  const results = searchResults[i];\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date and the multiplication answers in the practice sheet table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2023-08-06 Sunday'; New = '2023-08-07 Monday' },\n    @{ Old = '95\u00d714=1330'; New = '47\u00d763=2961' },\n    @{ Old = '63\u00d796=6048'; New = '22\u00d783=1826' },\n    @{ Old = '63\u00d778=4914'; New = '91\u00d746=4186' },\n    @{ Old = '49\u00d758=2842'; New = '72\u00d792=6624' },\n    @{ Old = '90\u00d778=7020'; New = '32\u00d732=1024' },\n    @{ Old = '43\u00d728=1204'; New = '38\u00d711=418' },\n    @{ Old = '47\u00d786=4042'; New = '96\u00d729=2784' },\n    @{ Old = '80\u00d774=5920'; New = '16\u00d733=528' },\n    @{ Old = '79\u00d742=3318'; New = '22\u00d766=1452' },\n    @{ Old = '33\u00d731=1023'; New = '76\u00d761=4636' },\n    @{ Old = '30\u00d777=2310'; New = '92\u00d718=1656' },\n    @{ Old = '65\u00d731=2015'; New = '75\u00d764=4800' },\n    @{ Old = '15\u00d716=240'; New = '14\u00d783=1162' },\n    @{ Old = '13\u00d739=507'; New = '85\u00d728=2380' },\n    @{ Old = '12\u00d733=396'; New = '41\u00d755=2255' },\n    @{ Old = '61\u00d753=3233'; New = '58\u00d717=986' },\n    @{ Old = '88\u00d720=1760'; New = '72\u00d765=4680' },\n    @{ Old = '65\u00d736=2340'; New = '94\u00d734=3196' },\n    @{ Old = '75\u00d771=5325'; New = '28\u00d783=2324' },\n    @{ Old = '92\u00d725=2300'; New = '34\u00d795=3230' },\n    @{ Old = '24\u00d732=768'; New = '94\u00d754=5076' },\n    @{ Old = '37\u00d796=3552'; New = '31\u00d754=1674' },\n    @{ Old = '74\u00d750=3700'; New = '94\u00d793=8742' },\n    @{ Old = '89\u00d779=7031'; New = '76\u00d729=2204' },\n    @{ Old = '44\u00d717=748'; New = '80\u00d722=1760' }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute(\n        [ref]$find.Text,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$true,\n        [ref]1,\n        [ref]$false,\n        [ref]$find.Replacement.Text,\n        [ref]2\n    ) | Out-Null\n}\n"}
